$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -1
$ws.Range("F3").Value = 3
$ws.Range("F5").Value = -1
$ws.Range("F6").Value = 0
$ws.Range("F11").Value = -8
$ws.Range("F13").Value = -1
$ws.Range("F16").Value = -6
$ws.Range("F20").Value = -5
$ws.Range("F22").Value = -1
$ws.Range("F25").Value = 3
$ws.Range("F26").Value = 6
$ws.Range("F27").Value = -6
